# Move existing files to cartesian sub folder
# - Fix D11/D21/D31/D41/D51 (previously buggy formulas referencing the
#   wrong row) to hold the plain static value -1.96, matching the other
#   "+90 degree" rows in the pattern.
# - Add a new column E with the euclidean distance
#   SQRT(B^2 + C^2) for every data row (2-61).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the mis-referencing formulas in column D -----------------------
$fixRows = @(11, 21, 31, 41, 51)
foreach ($r in $fixRows) {
    $ws.Range("D$r").Value = -1.96
}

# --- Add column E: distance = SQRT(B*B + C*C) ----------------------------
# First data row (E2) gets its own (non shared) formula.
$ws.Range("E2").Formula = "=SQRT(B2 * B2 + C2 * C2)"

# Remaining rows (E3:E61) are filled as one shared formula block.
$ws.Range("E3:E61").Formula = "=SQRT(B3 * B3 + C3 * C3)"

# --- Misc view/selection tweaks ------------------------------------------
[void]$ws.Range("D7").Select()
